$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'44.326.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.98%  '
$ws.Range("D3").Value = "'2.372.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'0.695"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.60%  '
$ws.Range("D6").Value = "'244.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.87%  '
$ws.Range("D7").Value = "'76.49"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +6.07%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = "'0.593"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +24.86%  '
$ws.Range("E10").Value = '  +4.97%  '
$ws.Range("D11").Value = "'58.03"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.07%  '
$ws.Range("D12").Value = "'32.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +17.11%  '
$ws.Range("D13").Value = "'7.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +18.11%  '
$ws.Range("D14").Value = "'0.109"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.57%  '
$ws.Range("D15").Value = "'2.725.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.25%  '
$ws.Range("D16").Value = "'17.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.53%  '
$ws.Range("D17").Value = "'0.926"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +7.63%  '
$ws.Range("D18").Value = "'2.371.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.61%  '
$ws.Range("D19").Value = "'44.335.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.96%  '
$ws.Range("E20").Value = '  +2.83%  '
$ws.Range("D21").Value = "'6.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.44%  '
$ws.Range("D22").Value = "'78.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.10%  '
$ws.Range("D23").Value = "'258.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.84%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("D25").Value = "'2.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.58%  '
$ws.Range("D26").Value = "'3.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.66%  '
$ws.Range("D27").Value = "'10.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.11%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = "'2.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.94%  '
$ws.Range("B29").Value = 'ImmutableX'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D29").Value = "'1.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +12.40%  '
$ws.Range("D30").Value = "'23.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.22%  '
$ws.Range("D31").Value = "'175.52"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.45%  '
$ws.Range("E32").Value = '  +0.33%  '
$ws.Range("E33").Value = '  +6.46%  '
$ws.Range("D34").Value = "'5.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.77%  '
$ws.Range("D35").Value = "'0.0762"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +9.75%  '
$ws.Range("E36").Value = '  +5.62%  '
$ws.Range("D37").Value = "'3.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.08%  '
$ws.Range("E38").Value = '  +1.12%  '
$ws.Range("D39").Value = "'6.60"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.05%  '
$ws.Range("D40").Value = "'0.0276"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.50%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = "'9.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.95%  '
$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D42").Value = "'19.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.41%  '
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("E44").Value = '  +15.84%  '
$ws.Range("D45").Value = "'1.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.55%  '
$ws.Range("D46").Value = "'0.102"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.93%  '
$ws.Range("D47").Value = "'1.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.55%  '
$ws.Range("D49").Value = "'102.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.10%  '
$ws.Range("D50").Value = "'4.47"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.14%  '
$ws.Range("D51").Value = "'1.470.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.39%  '
